# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-60) was re-sorted from
# descending chronological order (2003 .. 1607) to ascending chronological
# order (1607 .. 2003). The worker/doc/salary columns are constant for every
# row, so the only visible effect is that column E (Periodo Mora) and
# column F (Valor Mora) end up reversed relative to their original order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$valores = @(
    40000,40000,40000,40000,40000,40000,
    40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,
    40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,
    40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,
    40000,40000,38666
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
